$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "66.179.95"
$ws.Range("E2").Value2 = "  -0.16%  "
$ws.Range("D3").Value2 = "3.543.75"
$ws.Range("E3").Value2 = "  -0.24%  "
$ws.Range("E4").Value2 = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "602.67"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value2 = "  -0.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "145.69"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value2 = "  +0.97%  "
$ws.Range("D7").Value2 = "3.542.29"
$ws.Range("E7").Value2 = "  -0.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "1.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value2 = "  -0.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.494"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value2 = "  +1.02%  "
$ws.Range("E10").Value2 = "  -2.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "7.79"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value2 = "  -0.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.407"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value2 = "  -1.37%  "
$ws.Range("D13").Value2 = "4.143.92"
$ws.Range("E13").Value2 = "  -0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "0.0000201"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value2 = "  -2.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "28.94"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value2 = "  -4.00%  "
$ws.Range("D16").Value2 = "3.539.73"
$ws.Range("E16").Value2 = "  -0.48%  "
$ws.Range("E17").Value2 = "  +1.38%  "
$ws.Range("D18").Value2 = "66.103.40"
$ws.Range("E18").Value2 = "  -0.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "11.02"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value2 = "  -3.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "6.20"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value2 = "  +0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "14.64"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value2 = "  -1.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "416.54"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value2 = "  -3.40%  "
$ws.Range("E23").Value2 = "  -1.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "77.83"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value2 = "  -2.09%  "
$ws.Range("D25").Value2 = "3.684.30"
$ws.Range("E25").Value2 = "  -0.33%  "
$ws.Range("E26").Value2 = "  +0.13%  "
$ws.Range("E27").Value2 = "  -2.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "9.06"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value2 = "  -0.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "2.47"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value2 = "  -1.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "7.77"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value2 = "  -2.12%  "
$ws.Range("E31").Value2 = "  +0.11%  "
$ws.Range("D32").Value2 = "3.540.24"
$ws.Range("E32").Value2 = "  -0.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "0.156"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value2 = "  +2.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "24.35"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value2 = "  -4.15%  "
$ws.Range("E35").Value2 = "  +0.03%  "
$ws.Range("E37").Value2 = "  -10.37%  "
$ws.Range("B38").Value2 = "ImmutableX"
$ws.Range("C38").Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "1.60"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value2 = "  -7.15%  "
$ws.Range("B39").Value2 = "Monero"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "173.74"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value2 = "  -1.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "5.26"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value2 = "  -5.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.0819"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value2 = "  -3.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "5.07"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value2 = "  -2.42%  "
$ws.Range("E43").Value2 = "  -3.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "45.59"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value2 = "  -0.72%  "
$ws.Range("E45").Value2 = "  -6.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "0.998"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value2 = "  -0.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "2.41"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value2 = "  -4.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "7.08"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value2 = "  -0.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "22.69"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value2 = "  -2.33%  "
$ws.Range("E50").Value2 = "  -7.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "23.21"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value2 = "  -7.80%  "
